$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.008157135860865
$ws.Range("D2").Value = 1.027871511328864
$ws.Range("E2").Value = 1.010676782113318
$ws.Range("F2").Value = 1.022788911153612
$ws.Range("I2").Value = 1.029650528235211
$ws.Range("J2").Value = 1.013424596645301
$ws.Range("K2").Value = 1.030689977327416
$ws.Range("L2").Value = 1.013546252126332
$ws.Range("M2").Value = 1.025622259157619
$ws.Range("N2").Value = 1.008389760139752
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.009918853524015
$ws.Range("D3").Value = 1.028270653945178
$ws.Range("E3").Value = 1.012198607220225
$ws.Range("F3").Value = 1.024527811192356
$ws.Range("I3").Value = 1.029725471396968
$ws.Range("J3").Value = 1.014813763486936
$ws.Range("K3").Value = 1.030898040735165
$ws.Range("L3").Value = 1.014870194903636
$ws.Range("M3").Value = 1.027165357312983
$ws.Range("N3").Value = 1.008873429238299
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.011055844668206
$ws.Range("D4").Value = 1.028528447434878
$ws.Range("E4").Value = 1.01318088038086
$ws.Range("F4").Value = 1.025649382256573
$ws.Range("I4").Value = 1.029772114980503
$ws.Range("J4").Value = 1.015709535536505
$ws.Range("K4").Value = 1.031031454934845
$ws.Range("L4").Value = 1.015723955875829
$ws.Range("M4").Value = 1.028159807537883
$ws.Range("N4").Value = 1.009184666724216
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.011533144618985
$ws.Range("D5").Value = 1.028636708204138
$ws.Range("E5").Value = 1.013593255019767
$ws.Range("F5").Value = 1.026120040820188
$ws.Range("I5").Value = 1.029791281330322
$ws.Range("J5").Value = 1.016085386410795
$ws.Range("K5").Value = 1.031087250560783
$ws.Range("L5").Value = 1.016082190658627
$ws.Range("M5").Value = 1.028576921609032
$ws.Range("N5").Value = 1.00931510151637
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.011613245344475
$ws.Range("D6").Value = 1.028654878852546
$ws.Range("E6").Value = 1.013662461322285
$ws.Range("F6").Value = 1.026199017062729
$ws.Range("I6").Value = 1.029794473495223
$ws.Range("J6").Value = 1.01614845090259
$ws.Range("K6").Value = 1.031096601773843
$ws.Range("L6").Value = 1.016142299965408
$ws.Range("M6").Value = 1.028646901421496
$ws.Range("N6").Value = 1.0093369782093
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.011062225065976
$ws.Range("D7").Value = 1.028529894476086
$ws.Range("E7").Value = 1.013186392784138
$ws.Range("F7").Value = 1.025655674538123
$ws.Range("I7").Value = 1.029772372821
$ws.Range("J7").Value = 1.015714560531222
$ws.Range("K7").Value = 1.031032201625336
$ws.Range("L7").Value = 1.01572874530643
$ws.Range("M7").Value = 1.0281653847615
$ws.Range("N7").Value = 1.009186411203062
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.008753139356514
$ws.Range("D8").Value = 1.028006501881878
$ws.Range("E8").Value = 1.0111916053135
$ws.Range("F8").Value = 1.023377337950377
$ws.Range("I8").Value = 1.029676238734683
$ws.Range("J8").Value = 1.013894724427056
$ws.Range("K8").Value = 1.030760545020713
$ws.Range("L8").Value = 1.013994296571734
$ws.Range("M8").Value = 1.026144600721219
$ws.Range("N8").Value = 1.008553579004346
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.004660812932933
$ws.Range("D9").Value = 1.027080598561885
$ws.Range("E9").Value = 1.007657170216735
$ws.Range("F9").Value = 1.019334257081912
$ws.Range("I9").Value = 1.029492662180238
$ws.Range("J9").Value = 1.01066350910672
$ws.Range("K9").Value = 1.030272549766338
$ws.Range("L9").Value = 1.010915065350717
$ws.Range("M9").Value = 1.022552176623853
$ws.Range("N9").Value = 1.007425007589222
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.001915768734547
$ws.Range("D10").Value = 1.026460956289723
$ws.Range("E10").Value = 1.005287009748382
$ws.Range("F10").Value = 1.016618836919628
$ws.Range("I10").Value = 1.029360730725203
$ws.Range("J10").Value = 1.00849209672758
$ws.Range("K10").Value = 1.029940986572168
$ws.Range("L10").Value = 1.008846054973964
$ws.Range("M10").Value = 1.020135141993342
$ws.Range("N10").Value = 1.006663308168524
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.000722901755992
$ws.Range("D11").Value = 1.026192093312458
$ws.Range("E11").Value = 1.004257223544421
$ws.Range("F11").Value = 1.015438061907195
$ws.Range("I11").Value = 1.029301334300481
$ws.Range("J11").Value = 1.007547566427437
$ws.Range("K11").Value = 1.029795941338655
$ws.Range("L11").Value = 1.007946134561599
$ws.Range("M11").Value = 1.019083104073592
$ws.Range("N11").Value = 1.006331210227489
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.000279159213094
$ws.Range("D12").Value = 1.026092143331504
$ws.Range("E12").Value = 1.003874174065741
$ws.Range("F12").Value = 1.014998702464366
$ws.Range("I12").Value = 1.029278930456391
$ws.Range("J12").Value = 1.00719606425034
$ws.Range("K12").Value = 1.029741843520235
$ws.Range("L12").Value = 1.007611244016314
$ws.Range("M12").Value = 1.018691495021204
$ws.Range("N12").Value = 1.006207506627594
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.000374373637976
$ws.Range("D13").Value = 1.026113586651815
$ws.Range("E13").Value = 1.003956364229714
$ws.Range("F13").Value = 1.015092981570655
$ws.Range("I13").Value = 1.029283751611878
$ws.Range("J13").Value = 1.007271492863097
$ws.Range("K13").Value = 1.029753457708733
$ws.Range("L13").Value = 1.007683107480435
$ws.Range("M13").Value = 1.018775534570017
$ws.Range("N13").Value = 1.0062340572907
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.000686235416834
$ws.Range("D14").Value = 1.026183833089653
$ws.Range("E14").Value = 1.004225571688308
$ws.Range("F14").Value = 1.015401760075497
$ws.Range("I14").Value = 1.029299489358153
$ws.Range("J14").Value = 1.007518524739139
$ws.Range("K14").Value = 1.029791474114696
$ws.Range("L14").Value = 1.007918465142307
$ws.Range("M14").Value = 1.019050750687504
$ws.Range("N14").Value = 1.006320991973331
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.000878295845924
$ws.Range("D15").Value = 1.026227103359276
$ws.Range("E15").Value = 1.004391367098318
$ws.Range("F15").Value = 1.015591906497331
$ws.Range("I15").Value = 1.029309140655026
$ws.Range("J15").Value = 1.007670641073145
$ws.Range("K15").Value = 1.02981486792433
$ws.Range("L15").Value = 1.008063394144279
$ws.Range("M15").Value = 1.019220209215397
$ws.Range("N15").Value = 1.006374509068226
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.001994844176739
$ws.Range("D16").Value = 1.026478788249916
$ws.Range("E16").Value = 1.005355278307506
$ws.Range("F16").Value = 1.016697094541661
$ws.Range("I16").Value = 1.029364624823003
$ws.Range("J16").Value = 1.008554690221036
$ws.Range("K16").Value = 1.0299505816644
$ws.Range("L16").Value = 1.008905693620014
$ws.Range("M16").Value = 1.020204846087819
$ws.Range("N16").Value = 1.0066852999693
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.002694074330195
$ws.Range("D17").Value = 1.026636515912732
$ws.Range("E17").Value = 1.005958968073142
$ws.Range("F17").Value = 1.017389003265614
$ws.Range("I17").Value = 1.029398820732597
$ws.Range("J17").Value = 1.009108069659708
$ws.Range("K17").Value = 1.030035316254035
$ws.Range("L17").Value = 1.009432957420526
$ws.Range("M17").Value = 1.020821012719445
$ws.Range("N17").Value = 1.006879637371936
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.003101515319602
$ws.Range("D18").Value = 1.026728462299616
$ws.Range("E18").Value = 1.006310754388708
$ws.Range("F18").Value = 1.017792102688247
$ws.Range("I18").Value = 1.029418547711189
$ws.Range("J18").Value = 1.009430433292483
$ws.Range("K18").Value = 1.030084598063527
$ws.Range("L18").Value = 1.009740114137254
$ws.Range("M18").Value = 1.021179888002338
$ws.Range("N18").Value = 1.006992771765012
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.003240373570889
$ws.Range("D19").Value = 1.026759804540506
$ws.Range("E19").Value = 1.006430647986057
$ws.Range("F19").Value = 1.017929468454114
$ws.Range("I19").Value = 1.029425236981855
$ws.Range("J19").Value = 1.00954028136194
$ws.Range("K19").Value = 1.030101377718957
$ws.Range("L19").Value = 1.009844781400352
$ws.Range("M19").Value = 1.021302166785929
$ws.Range("N19").Value = 1.007031310634166
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.002619095917672
$ws.Range("D20").Value = 1.026619598753507
$ws.Range("E20").Value = 1.005894232717987
$ws.Range("F20").Value = 1.017314817712391
$ws.Range("I20").Value = 1.029395174481913
$ws.Range("J20").Value = 1.009048740143659
$ws.Range("K20").Value = 1.030026239766467
$ws.Range("L20").Value = 1.009376427179067
$ws.Range("M20").Value = 1.020754958224393
$ws.Range("N20").Value = 1.006858809533651
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.000594418246691
$ws.Range("D21").Value = 1.026163149529359
$ws.Range("E21").Value = 1.00414631187052
$ws.Range("F21").Value = 1.01531085382801
$ws.Range("I21").Value = 1.029294864411176
$ws.Range("J21").Value = 1.00744579840387
$ws.Range("K21").Value = 1.029780285345185
$ws.Range("L21").Value = 1.007849175416741
$ws.Range("M21").Value = 1.018969729598302
$ws.Range("N21").Value = 1.006295401517418
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 0.9993176018299114
$ws.Range("D22").Value = 1.02587568604674
$ws.Range("E22").Value = 1.003044186075239
$ws.Range("F22").Value = 1.014046434902394
$ws.Range("I22").Value = 1.029229820209068
$ws.Range("J22").Value = 1.006434129414482
$ws.Range("K22").Value = 1.029624362248677
$ws.Range("L22").Value = 1.006885336572081
$ws.Range("M22").Value = 1.01784244390925
$ws.Range("N22").Value = 1.005939150413415
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 0.9999948356264771
$ws.Range("D23").Value = 1.026028120679429
$ws.Range("E23").Value = 1.003628746888766
$ws.Range("F23").Value = 1.01471715519643
$ws.Range("I23").Value = 1.029264488757242
$ws.Range("J23").Value = 1.006970803409187
$ws.Range("K23").Value = 1.029707141463984
$ws.Range("L23").Value = 1.007396631702296
$ws.Range("M23").Value = 1.018440504087804
$ws.Range("N23").Value = 1.006128198647115
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.002652976695525
$ws.Range("D24").Value = 1.026627243055171
$ws.Range("E24").Value = 1.005923484875529
$ws.Range("F24").Value = 1.017348340451128
$ws.Range("I24").Value = 1.029396822742176
$ws.Range("J24").Value = 1.009075549874698
$ws.Range("K24").Value = 1.030030341480836
$ws.Range("L24").Value = 1.00940197195844
$ws.Range("M24").Value = 1.020784807027202
$ws.Range("N24").Value = 1.006868221414884
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.005721663359451
$ws.Range("D25").Value = 1.027320389354417
$ws.Range("E25").Value = 1.008573288179382
$ws.Range("F25").Value = 1.020382949875047
$ws.Range("I25").Value = 1.02954180262168
$ws.Range("J25").Value = 1.011501837092792
$ws.Range("K25").Value = 1.030399808295262
$ws.Range("L25").Value = 1.011713915783775
$ws.Range("M25").Value = 1.023484732094045
$ws.Range("N25").Value = 1.007718391891647
